$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.005723000769734084
$ws.Range("J2").Value = 0.005723000769734084
$ws.Range("S2").Value = 0.005723000769734084
$ws.Range("T2").Value = 0.005723000769734084

# Row 3 updates
$ws.Range("G3").Value = 30.199365
$ws.Range("H3").Value = 90.598095
$ws.Range("I3").Value = 0.6125398923302606
$ws.Range("J3").Value = 0.6125398923302606
$ws.Range("Q3").Value = 1.884289379175
$ws.Range("R3").Value = 16.958604412575
$ws.Range("S3").Value = 0.6125398923302606
$ws.Range("T3").Value = 0.6125398923302606

# Row 4 updates
$ws.Range("G4").Value = 18.820355
$ws.Range("H4").Value = 56.461065
$ws.Range("I4").Value = 0.3817371069000054
$ws.Range("J4").Value = 0.3817371069000054
$ws.Range("Q4").Value = 1.174296050225
$ws.Range("R4").Value = 10.568664452025
$ws.Range("S4").Value = 0.3817371069000054
$ws.Range("T4").Value = 0.3817371069000054
